$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3920
$ws.Range("C3").Value = 3692
$ws.Range("C4").Value = 2806
$ws.Range("C5").Value = 1947
$ws.Range("C6").Value = 1743
$ws.Range("C7").Value = 825
$ws.Range("C8").Value = 618
$ws.Range("C9").Value = 560
$ws.Range("C10").Value = 512
$ws.Range("B11").Value = "Kitchen & Dining"
$ws.Range("C11").Value = 505
